$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row/Col -> (Old, New) values straight from the table grid. We set the cell
# Range.Text directly (rather than Find.Execute/Replace) because this runtime's
# Find searches the whole document regardless of the Range it is invoked on,
# which would cause cross-cell collisions for duplicate equation text.
$replacements = @(
    @{Row=1; Col=1; Old="37×42=1554"; New="23×51=1173"},
    @{Row=1; Col=2; Old="85×14=1190"; New="24×13=312"},
    @{Row=1; Col=3; Old="96×48=4608"; New="85×74=6290"},
    @{Row=1; Col=4; Old="64×40=2560"; New="75×63=4725"},
    @{Row=1; Col=5; Old="77×43=3311"; New="69×66=4554"},
    @{Row=5; Col=1; Old="70×82=5740"; New="18×82=1476"},
    @{Row=5; Col=2; Old="40×59=2360"; New="59×95=5605"},
    @{Row=5; Col=3; Old="67×45=3015"; New="21×11=231"},
    @{Row=5; Col=4; Old="45×51=2295"; New="32×90=2880"},
    @{Row=5; Col=5; Old="78×17=1326"; New="90×61=5490"},
    @{Row=10; Col=1; Old="42×46=1932"; New="74×49=3626"},
    @{Row=10; Col=2; Old="83×39=3237"; New="11×75=825"},
    @{Row=10; Col=3; Old="48×52=2496"; New="62×70=4340"},
    @{Row=10; Col=4; Old="86×96=8256"; New="29×95=2755"},
    @{Row=10; Col=5; Old="39×12=468"; New="20×42=840"},
    @{Row=15; Col=1; Old="49×21=1029"; New="59×13=767"},
    @{Row=15; Col=2; Old="18×31=558"; New="37×81=2997"},
    @{Row=15; Col=3; Old="98×94=9212"; New="88×17=1496"},
    @{Row=15; Col=4; Old="21×16=336"; New="62×76=4712"},
    @{Row=15; Col=5; Old="27×94=2538"; New="70×56=3920"},
    @{Row=20; Col=1; Old="87×52=4524"; New="58×58=3364"},
    @{Row=20; Col=2; Old="11×75=825"; New="43×51=2193"},
    @{Row=20; Col=3; Old="20×39=780"; New="57×19=1083"},
    @{Row=20; Col=4; Old="96×59=5664"; New="60×33=1980"},
    @{Row=20; Col=5; Old="12×22=264"; New="78×36=2808"}
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $rng = $cell.Range
    # Cell.Range.Text carries trailing cell-mark control chars (CR + BEL);
    # strip them before sanity-checking against the expected old value.
    $before = $rng.Text.TrimEnd([char]13, [char]7)
    if ($before -ne $item.Old) {
        Write-Host "WARNING: cell ($($item.Row),$($item.Col)) expected `"$($item.Old)`" but found `"$before`""
    }
    $rng.Text = $item.New
}

Write-Host "All replacements applied."